$wb = $excel.ActiveWorkbook

# Rename sheet 6 ("constant_names") to "target"
$ws = $wb.Worksheets.Item(6)
$ws.Name = "target"

# Populate the new "target" sheet with merged constant/wavelength data
$ws.Range("A1").Value = "constant"
$ws.Range("B1").Value = "SB"
$ws.Range("A2").Value = "wavelength"
$ws.Range("B2").Value = 306
$ws.Range("C2").Value = 387

# Make "target" the active sheet/tab, with B1 selected
$ws.Activate()
$ws.Range("B1").Select()
